# Mengubah data toko ajat
# - Remove the computed "selisih" column (I): G - H formulas are dropped.
# - Insert a new header row at the top of the sheet with the field names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop column I (the G-H formula column) entirely.
$ws.Columns.Item(9).Delete()

# 2. Populate the (previously empty) header row (row 1) with the column names (written in
#    this particular order so the shared-string table is built in the same
#    sequence as the source data export).
$ws.Range("E1").Value = "harga_jual"
$ws.Range("A1").Value = "kategori"
$ws.Range("B1").Value = "nama_produk"
$ws.Range("C1").Value = "satuan"
$ws.Range("D1").Value = "harga_beli"
$ws.Range("F1").Value = "kecepatan_penjualan_per_hari"
$ws.Range("G1").Value = "jumlah_modal_stok"
$ws.Range("H1").Value = "jumlah_stok_sekarang"

# 3. Update the active selection to match the target workbook state.
$ws.Range("H1").Select()
